$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (matching source formatting)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.221.31"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.857.05"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "314.03"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "0.5091"
$ws.Range("E7").Value = "  -0.94%  "
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "0.08267"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "1.109"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "41.69"
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "6.201"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("D13").Value = "1.860.99"
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "20.23"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").Value = "7.188"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "0.00001099"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "90.99"
$ws.Range("E18").Value = "  -0.52%  "
$ws.Range("D19").Value = "0.06690"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "17.53"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "1.011"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").Value = "5.918"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").Value = "28.223.41"
$ws.Range("E23").Value = "  -0.24%  "
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").Value = "2.251"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "2.063.70"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("D27").Value = "159.83"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").Value = "20.59"
$ws.Range("E28").Value = "  -0.34%  "
$ws.Range("D29").Value = "2.365"
$ws.Range("E29").Value = "  -6.03%  "
$ws.Range("D30").Value = "126.04"
$ws.Range("E30").Value = "  +0.34%  "
$ws.Range("D31").Value = "0.1044"
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "1.023"
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "5.778"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("D34").Value = "3.626"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "0.02416"
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("D36").Value = "0.06431"
$ws.Range("E36").Value = "  -2.28%  "
$ws.Range("D37").Value = "9.040"
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("D38").Value = "0.2165"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").Value = "1.244"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.176"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.6385"
$ws.Range("E41").Value = "  -2.50%  "
$ws.Range("D42").Value = "4.918"
$ws.Range("E42").Value = "  -2.27%  "
$ws.Range("D43").Value = "11.06"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "0.5982"
$ws.Range("E44").Value = "  -2.57%  "
$ws.Range("D45").Value = "3.688"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "1.278"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "12.75"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("D48").Value = "1.969"
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("D49").Value = "1.201"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "120.59"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "0.06817"
$ws.Range("E51").Value = "  -1.24%  "
